$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 543-544, pushing the existing rows 543:589
# down to 545:591 (dimension grows from A1:R589 to A1:R591).
$ws.Range("543:544").Insert()

# --- New row 543 ---
$ws.Range("A543").Value = 8
$ws.Range("B543").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C543").Value = 'Coquimbo'
$ws.Range("D543").Value = 45223
$ws.Range("E543").Value = 4
$ws.Range("F543").Value = 100112021
$ws.Range("G543").Value = 'Ají'
$ws.Range("H543").Value = 'Inferno'
$ws.Range("I543").Value = 'Primera'
$ws.Range("J543").Value = 500
$ws.Range("K543").Value = 23000
$ws.Range("L543").Value = 24000
$ws.Range("M543").Value = 23500
$ws.Range("N543").Value = '$/caja 10 kilos'
$ws.Range("O543").Value = 'Región de Arica y Parinacota'
$ws.Range("P543").Value = 2350
$ws.Range("Q543").Value = 10
$ws.Range("R543").Value = 'Hortaliza'

# --- New row 544 ---
$ws.Range("A544").Value = 8
$ws.Range("B544").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C544").Value = 'Coquimbo'
$ws.Range("D544").Value = 45223
$ws.Range("E544").Value = 4
$ws.Range("F544").Value = 100112021
$ws.Range("G544").Value = 'Ají'
$ws.Range("H544").Value = 'Inferno'
$ws.Range("I544").Value = 'Segunda'
$ws.Range("J544").Value = 320
$ws.Range("K544").Value = 14000
$ws.Range("L544").Value = 15000
$ws.Range("M544").Value = 14500
$ws.Range("N544").Value = '$/caja 10 kilos'
$ws.Range("O544").Value = 'Región de Arica y Parinacota'
$ws.Range("P544").Value = 1450
$ws.Range("Q544").Value = 10
$ws.Range("R544").Value = 'Hortaliza'
